# Auto update Excel log
# Appends 4 new mmWave sensor log rows (rows 36-39) to the "mmWave" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mmWave")

$newRows = @(
    @("2026-02-01", "17:39:04", "17:00", "Living Room", "NO_MOTION_DETECTED", "Inactive"),
    @("2026-02-01", "17:39:14", "17:00", "Living Room", "PRESENCE_DETECTED",  "Active"),
    @("2026-02-01", "17:39:25", "17:00", "Living Room", "PRESENCE_DETECTED",  "Active"),
    @("2026-02-01", "17:39:35", "17:00", "Living Room", "PRESENCE_DETECTED",  "Active")
)

# Scratch cell, well outside the used range, used to stage date-like text
# (e.g. "2026-02-01") so it can be copied into the log as plain text without
# Excel auto-converting it to a date serial number, and without leaving any
# NumberFormat/style change behind on the destination cells.
$scratch = $ws.Cells.Item(200, 200)

$startRow = 36
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    for ($c = 1; $c -le 6; $c++) {
        $value = $row[$c - 1]
        $cell = $ws.Cells.Item($r, $c)
        if ($c -eq 1) {
            # Stage as text in the scratch cell, then copy/paste-values so the
            # destination keeps its original (default) formatting.
            $scratch.NumberFormat = "@"
            $scratch.Value = $value
            $scratch.Copy()
            $cell.PasteSpecial(-4163)  # xlPasteValues
        } else {
            $cell.Value = $value
        }
    }
}

$scratch.Clear()
$excel.CutCopyMode = $false
